$d = $word.ActiveDocument

function Get-ParaText($pp) {
    return $pp.Range.Text.TrimEnd([char]13)
}

function Find-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        if ((Get-ParaText $pp) -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParaTextNoHighlight($doc, $index, $newText) {
    $pp = $doc.Paragraphs.Item($index)
    $r = $pp.Range
    $sub = $doc.Range($r.Start, $r.End - 1)
    $sub.Text = $newText
    $pp2 = $doc.Paragraphs.Item($index)
    $pp2.Range.HighlightColorIndex = 0
}

function Delete-ParaByText($doc, $text) {
    $idx = Find-ParaIndexByText $doc $text
    if ($idx -ge 1) {
        $pp = $doc.Paragraphs.Item($idx)
        $pp.Range.Delete()
    }
}

# Step 1: remove the now-redundant paragraphs first (while their text is
# still unique), since the upcoming edits will create duplicate text.
Delete-ParaByText $d "Car moves to a free track position"
Delete-ParaByText $d "Car crashes to wall"
Delete-ParaByText $d "Car crashes to another running car"

# Step 2: rewrite the three yellow-highlighted "winner" paragraphs with the
# new equivalence-class descriptions, clearing their highlight.
$idx = Find-ParaIndexByText $d "Get winner when no car passed finish line"
Set-ParaTextNoHighlight $d $idx "Car moves to a free track position"

$idx = Find-ParaIndexByText $d "Get winner when one car passed finish line"
Set-ParaTextNoHighlight $d $idx "Car crashes into wall"

$idx = Find-ParaIndexByText $d "Get winner when two cars passed finish line"
Set-ParaTextNoHighlight $d $idx "Car crashes into another running car"

# Step 3: "Car crashes to another crashed car" -> "Car crashes into car
# that's on moving path"
$idx = Find-ParaIndexByText $d "Car crashes to another crashed car"
$pp = $d.Paragraphs.Item($idx)
$r = $pp.Range
$sub = $d.Range($r.Start, $r.End - 1)
$sub.Text = "Car crashes into car that's on moving path"

Write-Host "Final paragraph list:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Host "$i => [$(Get-ParaText $pp)]"
}
